# Insert a new row at position 147, shifting the existing rows 147-250 down
# to 148-251, then populate the newly inserted row 147 with the new weekly
# price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 147 (pushes old row147.. down by one)
$ws.Rows("147:147").Insert()

# Populate the new row 147 with data (same categorical values as the rest of
# the sheet, new date/price figures for this week's observation)
$ws.Range("A147").Value = 10
$ws.Range("B147").Value = 'Vega Modelo de Temuco'
$ws.Range("C147").Value = 'La Araucanía'
$ws.Range("D147").Value = 44762
$ws.Range("E147").Value = 9
$ws.Range("F147").Value = 100112052
$ws.Range("G147").Value = 'Albahaca'
$ws.Range("H147").Value = 'Sin especificar'
$ws.Range("I147").Value = 'Primera'
$ws.Range("J147").Value = 45
$ws.Range("K147").Value = 5500
$ws.Range("L147").Value = 5500
$ws.Range("M147").Value = 5500
$ws.Range("N147").Value = '$/paquete'
$ws.Range("O147").Value = 'Región de Arica y Parinacota'
$ws.Range("P147").Value = 5500
$ws.Range("Q147").Value = 1
$ws.Range("R147").Value = 'Hortaliza'

# Make sure the date cell keeps the date/time number format used elsewhere
# in column D.
$ws.Range("D147").NumberFormat = $ws.Range("D148").NumberFormat
